$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.689.19'
$ws.Range("E2").Value = '  -3.55%  '

$ws.Range("D3").Value = '1.745.98'
$ws.Range("E3").Value = '  -5.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4926'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.41%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.55'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -21.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05958'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -12.35%  '

$ws.Range("D11").Value = '1.744.96'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06772'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -13.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -23.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.468'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -11.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -12.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5798'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -26.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '25.756.46'
$ws.Range("E19").Value = '  -3.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -17.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006526'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -17.86%  '

$ws.Range("D22").Value = '1.962.41'
$ws.Range("E22").Value = '  -5.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.988'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -13.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.063'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -15.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.930'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -15.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.493'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.839'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -17.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -14.51%  '

$ws.Range("E30").Value = '  -8.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.766'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08087'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.355'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -18.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04436'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.0000'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.669'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.013'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6045'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -17.61%  '

$ws.Range("E39").Value = '  -12.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.036'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.81%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.12%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01494'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -13.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7748'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -14.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.208'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -11.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.3748'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -22.32%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05119'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.29%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1080'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -13.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.940'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -23.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.23'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -13.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -12.38%  '
